# Daily update at 8 AM UTC
# Appends the next day's data row (row 57) to the "Wins Over Time" sheet
# and moves the "last row" date style from the old last row (56) to the
# new last row (57).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 56 was previously the last row and used the "date only, last row"
# number format (YYYY-MM-DD). Since it is no longer the last row, it
# reverts to the regular date/time number format used by all the other
# data rows.
$ws.Range("A56").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add the new data row (57) for 2025-05-20.
$ws.Range("A57").Value = 45797
$ws.Range("A57").NumberFormat = "YYYY-MM-DD"
$ws.Range("B57").Value = 238
$ws.Range("C57").Value = 248
$ws.Range("D57").Value = 238
